$wb = $excel.ActiveWorkbook

# --- Outputs sheet ---------------------------------------------------------
$wsOutputs = $wb.Worksheets.Item("Outputs")
# "winch lock" -> "Winch Lock"
$wsOutputs.Range("E14").Value = "Winch Lock"

# --- Other Inputs sheet -----------------------------------------------------
$wsOther = $wb.Worksheets.Item("Other Inputs")

# Row 11: Cyprus digital 1 - "hold to fire winch"
$wsOther.Range("A11").Value = "CYPRUS STUFF"
$wsOther.Range("B11").Value = "Digital"
$wsOther.Range("C11").Value = "D"
$wsOther.Range("D11").Value = 1
$wsOther.Range("E11").Value = "Hold to fire winch"

# Row 12: Cyprus digital 2 - "push to lock / unlock"
$wsOther.Range("C12").Value = "D"
$wsOther.Range("D12").Value = 2
$wsOther.Range("E12").Value = "Push to lock / unlock"

# Row 13: Cyprus digital 3 (unused)
$wsOther.Range("C13").Value = "D"
$wsOther.Range("D13").Value = 3

# Row 14: Cyprus digital 4 (unused)
$wsOther.Range("C14").Value = "D"
$wsOther.Range("D14").Value = 4

# Row 15: Cyprus digital 5 (unused)
$wsOther.Range("C15").Value = "D"
$wsOther.Range("D15").Value = 5

# Row 16: Cyprus digital 6 (unused)
$wsOther.Range("C16").Value = "D"
$wsOther.Range("D16").Value = 6

# Row 17: Cyprus analog 1 - "change winch speed"
$wsOther.Range("B17").Value = "Dial"
$wsOther.Range("C17").Value = "A"
$wsOther.Range("D17").Value = 1
$wsOther.Range("E17").Value = "Change winch speed"

# --- Restore the cursor positions recorded in the author's edit ------------
$wsOther.Activate() | Out-Null
$wsOther.Range("B18").Select() | Out-Null
$wsOutputs.Activate() | Out-Null
$wsOutputs.Range("E15").Select() | Out-Null
$wsOther.Activate() | Out-Null
